$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Rename Sheet1 -> Credentials
$ws1.Name = "Credentials"

# Update credential rows on the (renamed) Credentials sheet.
# Set B1 before A1 so new shared-strings are appended in the same order
# the source workbook used (Siri@1234 before sirisha@regal-us.com).
$ws1.Range("B1").Value = "Siri@1234"
$ws1.Range("A1").Value = "sirisha@regal-us.com"
$ws1.Range("A2").Value = "Chethan@regal-us.com"
$ws1.Range("B2").Value = "Regal@123"

# Clear out the leftover credential row on Sheet2 without disturbing which
# sheet/tab is active (Credentials stays the active tab).
[void]$ws2.Range("A1:B1").Select()
$ws1.Activate()
$ws2.Range("A1:B1").ClearContents()
